$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.144.06'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.783.97'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '''226.16'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').Value = '''0.548'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').Value = '''0.0687'
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('D12').Value = '2.040.35'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '''10.98'
$ws.Range('E13').Value = '  -4.85%  '
$ws.Range('D14').Value = '1.770.98'
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = '34.128.29'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = '''0.622'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '''67.63'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').Value = '''245.54'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').Value = '''2.05'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').Value = '''162.14'
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').Value = '''7.14'
$ws.Range('E26').Value = '  -0.30%  '
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('D29').Value = '''1.01'
$ws.Range('E29').Value = '  +0.33%  '
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('D31').Value = '''1.23'
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('D32').Value = '''3.72'
$ws.Range('E32').Value = '  +1.56%  '
$ws.Range('D33').Value = '''3.74'
$ws.Range('E33').Value = '  +3.26%  '
$ws.Range('D34').Value = '''1.79'
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('D35').Value = '1.447.73'
$ws.Range('E35').Value = '  +3.04%  '
$ws.Range('D36').Value = '''2.52'
$ws.Range('E36').Value = '  +8.10%  '
$ws.Range('D37').Value = '''0.653'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('E38').Value = '  +1.23%  '
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').Value = '''82.59'
$ws.Range('E40').Value = '  +3.26%  '
$ws.Range('D42').Value = '''2.72'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('D43').Value = '''0.915'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('D44').Value = '''13.69'
$ws.Range('E44').Value = '  +2.31%  '
$ws.Range('D45').Value = '''0.0519'
$ws.Range('E45').Value = '  +2.40%  '
$ws.Range('D46').Value = '''6.07'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').Value = '1.940.19'
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0132'
$ws.Range('E49').Value = '  -5.49%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '''104.84'
$ws.Range('E50').Value = '  -1.94%  '

# Reset style on text-forced numeric-looking cells back to default (remove quote-prefix style)
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
